$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = -7.895
$ws.Range("C7").Value = -12.844
$ws.Range("A8").Value = -22.188
$ws.Range("A10").Value = -21.777
$ws.Range("A12").Value = -21.628
$ws.Range("C15").Value = -13.728
$ws.Range("A18").Value = -21.841
$ws.Range("C18").Value = -11.773
$ws.Range("D18").Value = -7.738000000000001
$ws.Range("D19").Value = -8.113
$ws.Range("C20").Value = -12.673
$ws.Range("D27").Value = -8.535
$ws.Range("C29").Value = -12.181
$ws.Range("C30").Value = -12.569
$ws.Range("C31").Value = -13.363
$ws.Range("D31").Value = -8.425999999999998
$ws.Range("A37").Value = -20.029
$ws.Range("D38").Value = -7.726999999999999
$ws.Range("C40").Value = -12.782
$ws.Range("D42").Value = -8.257999999999999
$ws.Range("D44").Value = -7.316000000000001
$ws.Range("D47").Value = -7.394999999999999
$ws.Range("C50").Value = -13.326
$ws.Range("A55").Value = -21.894
$ws.Range("D58").Value = -8.400999999999998
$ws.Range("D65").Value = -7.437
$ws.Range("A68").Value = -21.593
$ws.Range("D73").Value = -8.105
$ws.Range("C76").Value = -13.045
$ws.Range("A77").Value = -20.843
$ws.Range("A78").Value = -19.951
$ws.Range("A81").Value = -21.836
$ws.Range("A82").Value = -22.152
$ws.Range("C87").Value = -13.193
$ws.Range("C88").Value = -13.085
$ws.Range("D90").Value = -7.452
$ws.Range("D94").Value = -7.308000000000002
$ws.Range("D95").Value = -7.718999999999999
$ws.Range("C96").Value = -12.665
$ws.Range("C98").Value = -13.23
$ws.Range("C101").Value = -13.048
$ws.Range("D101").Value = -8.061999999999999
$ws.Range("C102").Value = -13.086
